$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("06_API一覧")

$ws.Cells.Item(32, 3).Value = "/api/v1/instructor/courses/{courseId}/publish"
$ws.Cells.Item(32, 4).Value = "POST"
$ws.Cells.Item(32, 5).Value = "コース公開（講師）"
$ws.Cells.Item(32, 6).Value = "app"
$ws.Cells.Item(32, 7).Value = "instructor_owner"
$ws.Cells.Item(32, 8).Value = "AUDIT_LOG"
$ws.Cells.Item(32, 9).Value = "-"
$ws.Cells.Item(32, 10).Value = "CourseDetailView"
$ws.Cells.Item(32, 11).Value = "200, 403, 423"
$ws.Cells.Item(32, 12).Value = "v1.2"
$ws.Cells.Item(32, 13).Value = "ownerUserId一致の講師がコースを公開。statusをactiveに変更。operator作成コースの場合、CourseMember.role=instructor→instructor_ownerへ昇格（委譲完了）。"
$ws.Cells.Item(32, 2).Value = "API-031"
$ws.Cells.Item(33, 3).Value = "/api/v1/instructor/analytics/overview"
$ws.Cells.Item(33, 4).Value = "GET"
$ws.Cells.Item(33, 5).Value = "売上分析サマリ"
$ws.Cells.Item(33, 6).Value = "app"
$ws.Cells.Item(33, 7).Value = "instructor"
$ws.Cells.Item(33, 8).Value = "-"
$ws.Cells.Item(33, 9).Value = "-"
$ws.Cells.Item(33, 10).Value = "GenericListResponse"
$ws.Cells.Item(33, 11).NumberFormat = "@"
$ws.Cells.Item(33, 11).Value = "200"
$ws.Cells.Item(33, 12).Value = "KEEP"
$ws.Cells.Item(33, 13).Value = "期間別の収益、成約率、受講生離脱ポイント等の分析サマリーを取得。"
$ws.Cells.Item(33, 2).Value = "API-032"
$ws.Cells.Item(34, 3).Value = "/api/v1/instructor/settings/payout"
$ws.Cells.Item(34, 4).Value = "PUT"
$ws.Cells.Item(34, 5).Value = "振込先銀行設定"
$ws.Cells.Item(34, 6).Value = "app"
$ws.Cells.Item(34, 7).Value = "instructor_owner"
$ws.Cells.Item(34, 8).Value = "-"
$ws.Cells.Item(34, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(34, 10).Value = "SuccessResponse"
$ws.Cells.Item(34, 11).NumberFormat = "@"
$ws.Cells.Item(34, 11).Value = "200"
$ws.Cells.Item(34, 12).Value = "KEEP"
$ws.Cells.Item(34, 13).Value = "講師への報酬振込先口座情報をStripe Connect等と連携して設定・更新。"
$ws.Cells.Item(34, 2).Value = "API-033"
$ws.Cells.Item(35, 3).Value = "/api/v1/instructor/courses/{courseId}/syllabus"
$ws.Cells.Item(35, 4).Value = "GET"
$ws.Cells.Item(35, 5).Value = "シラバス構造取得"
$ws.Cells.Item(35, 6).Value = "app"
$ws.Cells.Item(35, 7).Value = "instructor"
$ws.Cells.Item(35, 8).Value = "-"
$ws.Cells.Item(35, 9).Value = "-"
$ws.Cells.Item(35, 10).Value = "CourseDetailView"
$ws.Cells.Item(35, 11).NumberFormat = "@"
$ws.Cells.Item(35, 11).Value = "200"
$ws.Cells.Item(35, 12).Value = "KEEP"
$ws.Cells.Item(35, 13).Value = "編集中の章立てとレッスン構成をツリー形式で取得（image_32dd18の構成用）。"
$ws.Cells.Item(35, 2).Value = "API-034"
$ws.Cells.Item(36, 3).Value = "/api/v1/instructor/courses/{courseId}/sections"
$ws.Cells.Item(36, 4).Value = "POST"
$ws.Cells.Item(36, 5).Value = "セクション追加"
$ws.Cells.Item(36, 6).Value = "app"
$ws.Cells.Item(36, 7).Value = "instructor"
$ws.Cells.Item(36, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(36, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(36, 10).Value = "CourseDetailView"
$ws.Cells.Item(36, 11).Value = "201, 423"
$ws.Cells.Item(36, 12).Value = "KEEP"
$ws.Cells.Item(36, 13).Value = "カリキュラム内に新しい「章」を追加。423チェック対象。"
$ws.Cells.Item(36, 2).Value = "API-035"
$ws.Cells.Item(37, 3).Value = "/api/v1/instructor/sections/{sectionId}"
$ws.Cells.Item(37, 4).Value = "PUT"
$ws.Cells.Item(37, 5).Value = "セクション編集"
$ws.Cells.Item(37, 6).Value = "app"
$ws.Cells.Item(37, 7).Value = "instructor"
$ws.Cells.Item(37, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(37, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(37, 10).Value = "SuccessResponse"
$ws.Cells.Item(37, 11).Value = "200, 423"
$ws.Cells.Item(37, 12).Value = "KEEP"
$ws.Cells.Item(37, 13).Value = "セクション名の変更、表示順の並び替え、および削除。"
$ws.Cells.Item(37, 2).Value = "API-036"
$ws.Cells.Item(38, 3).Value = "/api/v1/instructor/sections/{sectionId}"
$ws.Cells.Item(38, 4).Value = "DELETE"
$ws.Cells.Item(38, 5).Value = "セクション削除"
$ws.Cells.Item(38, 6).Value = "app"
$ws.Cells.Item(38, 7).Value = "instructor"
$ws.Cells.Item(38, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(38, 9).Value = "-"
$ws.Cells.Item(38, 10).Value = "SuccessResponse"
$ws.Cells.Item(38, 11).Value = "200, 423"
$ws.Cells.Item(38, 12).Value = "KEEP"
$ws.Cells.Item(38, 13).Value = "指定されたセクションと、配下のレッスン紐付けを削除（論理/物理）。"
$ws.Cells.Item(38, 2).Value = "API-037"
$ws.Cells.Item(39, 3).Value = "/api/v1/instructor/sections/{sectionId}/lessons"
$ws.Cells.Item(39, 4).Value = "POST"
$ws.Cells.Item(39, 5).Value = "レッスン作成"
$ws.Cells.Item(39, 6).Value = "app"
$ws.Cells.Item(39, 7).Value = "instructor"
$ws.Cells.Item(39, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(39, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(39, 10).Value = "SuccessResponse"
$ws.Cells.Item(39, 11).Value = "201, 423"
$ws.Cells.Item(39, 12).Value = "KEEP"
$ws.Cells.Item(39, 13).Value = "指定セクションに新規レッスンを追加（座学/動画/課題/ライブ）（image_32dcf8反映）。"
$ws.Cells.Item(39, 2).Value = "API-038"
$ws.Cells.Item(40, 3).Value = "/api/v1/instructor/lessons/{lessonId}"
$ws.Cells.Item(40, 4).Value = "GET"
$ws.Cells.Item(40, 5).Value = "レッスン詳細取得"
$ws.Cells.Item(40, 6).Value = "app"
$ws.Cells.Item(40, 7).Value = "instructor"
$ws.Cells.Item(40, 8).Value = "-"
$ws.Cells.Item(40, 9).Value = "-"
$ws.Cells.Item(40, 10).Value = "GenericDetailView"
$ws.Cells.Item(40, 11).NumberFormat = "@"
$ws.Cells.Item(40, 11).Value = "200"
$ws.Cells.Item(40, 12).Value = "KEEP"
$ws.Cells.Item(40, 13).Value = "編集画面用。動画URL、リッチテキスト本文、課題要件、Drip設定の詳細。"
$ws.Cells.Item(40, 2).Value = "API-039"
$ws.Cells.Item(41, 3).Value = "/api/v1/instructor/lessons/{lessonId}"
$ws.Cells.Item(41, 4).Value = "PUT"
$ws.Cells.Item(41, 5).Value = "レッスン編集"
$ws.Cells.Item(41, 6).Value = "app"
$ws.Cells.Item(41, 7).Value = "instructor"
$ws.Cells.Item(41, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(41, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(41, 10).Value = "SuccessResponse"
$ws.Cells.Item(41, 11).Value = "200, 423"
$ws.Cells.Item(41, 12).Value = "KEEP"
$ws.Cells.Item(41, 13).Value = "コンテンツの修正、Drip(解禁日)、先行ロック条件の保存。423チェック対象。"
$ws.Cells.Item(41, 2).Value = "API-040"
$ws.Cells.Item(42, 3).Value = "/api/v1/instructor/lessons/{lessonId}"
$ws.Cells.Item(42, 4).Value = "DELETE"
$ws.Cells.Item(42, 5).Value = "レッスン削除"
$ws.Cells.Item(42, 6).Value = "app"
$ws.Cells.Item(42, 7).Value = "instructor"
$ws.Cells.Item(42, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(42, 9).Value = "-"
$ws.Cells.Item(42, 10).Value = "SuccessResponse"
$ws.Cells.Item(42, 11).Value = "200, 423"
$ws.Cells.Item(42, 12).Value = "KEEP"
$ws.Cells.Item(42, 13).Value = "レッスンをカリキュラムから削除。"
$ws.Cells.Item(42, 2).Value = "API-041"
$ws.Cells.Item(43, 3).Value = "/api/v1/instructor/lessons/{lessonId}/live"
$ws.Cells.Item(43, 4).Value = "PUT"
$ws.Cells.Item(43, 5).Value = "ライブ配信設定"
$ws.Cells.Item(43, 6).Value = "app"
$ws.Cells.Item(43, 7).Value = "instructor"
$ws.Cells.Item(43, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(43, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(43, 10).Value = "SuccessResponse"
$ws.Cells.Item(43, 11).Value = "200, 423"
$ws.Cells.Item(43, 12).Value = "KEEP"
$ws.Cells.Item(43, 13).Value = "ZoomやYouTube LiveのURL、および配信日時の設定（カレンダー連動）。"
$ws.Cells.Item(43, 2).Value = "API-042"
$ws.Cells.Item(44, 3).Value = "/api/v1/instructor/courses/{courseId}/submissions"
$ws.Cells.Item(44, 4).Value = "GET"
$ws.Cells.Item(44, 5).Value = "提出一覧取得"
$ws.Cells.Item(44, 6).Value = "app"
$ws.Cells.Item(44, 7).Value = "instructor"
$ws.Cells.Item(44, 8).Value = "-"
$ws.Cells.Item(44, 9).Value = "-"
$ws.Cells.Item(44, 10).Value = "CourseDetailView"
$ws.Cells.Item(44, 11).NumberFormat = "@"
$ws.Cells.Item(44, 11).Value = "200"
$ws.Cells.Item(44, 12).Value = "KEEP"
$ws.Cells.Item(44, 13).Value = "受講者別の課題進捗、未確認提出物、最終提出日時の一覧（image_32dcd7反映）。"
$ws.Cells.Item(44, 2).Value = "API-043"
$ws.Cells.Item(45, 3).Value = "/api/v1/instructor/submissions/{submissionId}/evaluation"
$ws.Cells.Item(45, 4).Value = "PATCH"
$ws.Cells.Item(45, 5).Value = "評価・採点実行"
$ws.Cells.Item(45, 6).Value = "app"
$ws.Cells.Item(45, 7).Value = "instructor, assistant"
$ws.Cells.Item(45, 8).Value = "THREAD_REPLY(AUTO)"
$ws.Cells.Item(45, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(45, 10).Value = "SubmissionView"
$ws.Cells.Item(45, 11).NumberFormat = "@"
$ws.Cells.Item(45, 11).Value = "200"
$ws.Cells.Item(45, 12).Value = "KEEP"
$ws.Cells.Item(45, 13).Value = "提出物への合否、講評を送信。合格時は自動で次の教材を解禁。フィードバックはスレッドへ。"
$ws.Cells.Item(45, 2).Value = "API-044"
$ws.Cells.Item(46, 3).Value = "/api/v1/instructor/courses/{courseId}/members"
$ws.Cells.Item(46, 4).Value = "GET"
$ws.Cells.Item(46, 5).Value = "受講者名簿取得"
$ws.Cells.Item(46, 6).Value = "app"
$ws.Cells.Item(46, 7).Value = "instructor"
$ws.Cells.Item(46, 8).Value = "-"
$ws.Cells.Item(46, 9).Value = "-"
$ws.Cells.Item(46, 10).Value = "CourseDetailView"
$ws.Cells.Item(46, 11).NumberFormat = "@"
$ws.Cells.Item(46, 11).Value = "200"
$ws.Cells.Item(46, 12).Value = "KEEP"
$ws.Cells.Item(46, 13).Value = "全参加ユーザーの属性、現在の進捗状況、メールアドレス等の名簿取得。"
$ws.Cells.Item(46, 2).Value = "API-045"
$ws.Cells.Item(47, 3).Value = "/api/v1/instructor/courses/{courseId}/members/{userId}/role"
$ws.Cells.Item(47, 4).Value = "PATCH"
$ws.Cells.Item(47, 5).Value = "講座内ロール変更"
$ws.Cells.Item(47, 6).Value = "app"
$ws.Cells.Item(47, 7).Value = "instructor_owner"
$ws.Cells.Item(47, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(47, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(47, 10).Value = "CourseDetailView"
$ws.Cells.Item(47, 11).Value = "200, 423"
$ws.Cells.Item(47, 12).Value = "KEEP"
$ws.Cells.Item(47, 13).Value = "特定ユーザーの役割（講師、アシスタント、受講生）を動的に切り替え。"
$ws.Cells.Item(47, 2).Value = "API-046"
$ws.Cells.Item(48, 3).Value = "/api/v1/instructor/courses/{courseId}/members/{userId}/revoke"
$ws.Cells.Item(48, 4).Value = "POST"
$ws.Cells.Item(48, 5).Value = "受講権限剥奪"
$ws.Cells.Item(48, 6).Value = "app"
$ws.Cells.Item(48, 7).Value = "instructor_owner"
$ws.Cells.Item(48, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(48, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(48, 10).Value = "CourseDetailView"
$ws.Cells.Item(48, 11).Value = "201, 423"
$ws.Cells.Item(48, 12).Value = "KEEP"
$ws.Cells.Item(48, 13).Value = "特定ユーザーの受講権限を剥奪（revoked）し、コンテンツへのアクセスを遮断。"
$ws.Cells.Item(48, 2).Value = "API-047"
$ws.Cells.Item(49, 3).Value = "/api/v1/instructor/courses/{courseId}/export"
$ws.Cells.Item(49, 4).Value = "GET"
$ws.Cells.Item(49, 5).Value = "受講者CSV出力"
$ws.Cells.Item(49, 6).Value = "app"
$ws.Cells.Item(49, 7).Value = "instructor_owner"
$ws.Cells.Item(49, 8).Value = "-"
$ws.Cells.Item(49, 9).Value = "-"
$ws.Cells.Item(49, 10).Value = "CourseDetailView"
$ws.Cells.Item(49, 11).NumberFormat = "@"
$ws.Cells.Item(49, 11).Value = "200"
$ws.Cells.Item(49, 12).Value = "KEEP"
$ws.Cells.Item(49, 13).Value = "受講者名簿と学習進捗データをCSV形式でエクスポートする。"
$ws.Cells.Item(49, 2).Value = "API-048"
$ws.Cells.Item(50, 3).Value = "/api/v1/courses/{courseId}/channels"
$ws.Cells.Item(50, 4).Value = "GET"
$ws.Cells.Item(50, 5).Value = "チャンネル一覧取得"
$ws.Cells.Item(50, 6).Value = "app"
$ws.Cells.Item(50, 7).Value = "all_in_course"
$ws.Cells.Item(50, 8).Value = "-"
$ws.Cells.Item(50, 9).Value = "-"
$ws.Cells.Item(50, 10).Value = "CourseChannelListResponse"
$ws.Cells.Item(50, 11).NumberFormat = "@"
$ws.Cells.Item(50, 11).Value = "200"
$ws.Cells.Item(50, 12).Value = "KEEP"
$ws.Cells.Item(50, 13).Value = "サイドメニュー構築用。general, announcement等のチャンネル種別を含む。"
$ws.Cells.Item(50, 2).Value = "API-049"
$ws.Cells.Item(51, 3).Value = "/api/v1/courses/{courseId}/channels"
$ws.Cells.Item(51, 4).Value = "POST"
$ws.Cells.Item(51, 5).Value = "チャンネル作成"
$ws.Cells.Item(51, 6).Value = "app"
$ws.Cells.Item(51, 7).Value = "instructor_owner"
$ws.Cells.Item(51, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(51, 9).Value = "CourseChannelCreateRequest"
$ws.Cells.Item(51, 10).Value = "CourseDetailView"
$ws.Cells.Item(51, 11).Value = "201, 423"
$ws.Cells.Item(51, 12).Value = "KEEP"
$ws.Cells.Item(51, 13).Value = "コース内に新しいカスタムチャンネルを追加。"
$ws.Cells.Item(51, 2).Value = "API-050"
$ws.Cells.Item(52, 3).Value = "/api/v1/channels/{channelId}"
$ws.Cells.Item(52, 4).Value = "PUT"
$ws.Cells.Item(52, 5).Value = "チャンネル編集"
$ws.Cells.Item(52, 6).Value = "app"
$ws.Cells.Item(52, 7).Value = "instructor_owner"
$ws.Cells.Item(52, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(52, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(52, 10).Value = "CourseChannelDetailView"
$ws.Cells.Item(52, 11).Value = "200, 423"
$ws.Cells.Item(52, 12).Value = "KEEP"
$ws.Cells.Item(52, 13).Value = "チャンネル名、説明文、閲覧制限、アーカイブ状態の更新。"
$ws.Cells.Item(52, 2).Value = "API-051"
$ws.Cells.Item(53, 3).Value = "/api/v1/channels/{channelId}"
$ws.Cells.Item(53, 4).Value = "DELETE"
$ws.Cells.Item(53, 5).Value = "チャンネル削除"
$ws.Cells.Item(53, 6).Value = "app"
$ws.Cells.Item(53, 7).Value = "instructor_owner"
$ws.Cells.Item(53, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(53, 9).Value = "-"
$ws.Cells.Item(53, 10).Value = "SuccessResponse"
$ws.Cells.Item(53, 11).Value = "200, 423"
$ws.Cells.Item(53, 12).Value = "KEEP"
$ws.Cells.Item(53, 13).Value = "チャンネルを論理削除。過去ログは監査用に保持。"
$ws.Cells.Item(53, 2).Value = "API-052"
$ws.Cells.Item(54, 3).Value = "/api/v1/channels/{channelId}/messages"
$ws.Cells.Item(54, 4).Value = "GET"
$ws.Cells.Item(54, 5).Value = "メッセージ履歴取得"
$ws.Cells.Item(54, 6).Value = "app"
$ws.Cells.Item(54, 7).Value = "all_in_course"
$ws.Cells.Item(54, 8).Value = "閲覧のみ可"
$ws.Cells.Item(54, 9).Value = "-"
$ws.Cells.Item(54, 10).Value = "MessageListResponse"
$ws.Cells.Item(54, 11).NumberFormat = "@"
$ws.Cells.Item(54, 11).Value = "200"
$ws.Cells.Item(54, 12).Value = "KEEP"
$ws.Cells.Item(54, 13).Value = "スレッド親（ルート）メッセージの一覧取得。課題の相互閲覧可（image_32d9b1反映）。"
$ws.Cells.Item(54, 2).Value = "API-053"
$ws.Cells.Item(55, 3).Value = "/api/v1/channels/{channelId}/messages"
$ws.Cells.Item(55, 4).Value = "POST"
$ws.Cells.Item(55, 5).Value = "メッセージ投稿"
$ws.Cells.Item(55, 6).Value = "app"
$ws.Cells.Item(55, 7).Value = "all_in_course"
$ws.Cells.Item(55, 8).Value = "threads_only(AUTO)"
$ws.Cells.Item(55, 9).Value = "CourseMessageCreateRequest"
$ws.Cells.Item(55, 10).Value = "CourseChannelDetailView"
$ws.Cells.Item(55, 11).NumberFormat = "@"
$ws.Cells.Item(55, 11).Value = "201"
$ws.Cells.Item(55, 12).Value = "KEEP"
$ws.Cells.Item(55, 13).Value = "チャンネルへの新規投稿。親メッセージとして保存。threads_onlyを強制適用。"
$ws.Cells.Item(55, 2).Value = "API-054"
$ws.Cells.Item(56, 3).Value = "/api/v1/messages/{messageId}/thread"
$ws.Cells.Item(56, 4).Value = "GET"
$ws.Cells.Item(56, 5).Value = "スレッド詳細取得"
$ws.Cells.Item(56, 6).Value = "app"
$ws.Cells.Item(56, 7).Value = "all_in_course"
$ws.Cells.Item(56, 8).Value = "-"
$ws.Cells.Item(56, 9).Value = "-"
$ws.Cells.Item(56, 10).Value = "GenericListResponse"
$ws.Cells.Item(56, 11).NumberFormat = "@"
$ws.Cells.Item(56, 11).Value = "200"
$ws.Cells.Item(56, 12).Value = "KEEP"
$ws.Cells.Item(56, 13).Value = "特定メッセージに紐づく返信メッセージ一覧（スレッドビュー）を全件取得。"
$ws.Cells.Item(56, 2).Value = "API-055"
$ws.Cells.Item(57, 3).Value = "/api/v1/messages/{messageId}/replies"
$ws.Cells.Item(57, 4).Value = "POST"
$ws.Cells.Item(57, 5).Value = "スレッド返信投稿"
$ws.Cells.Item(57, 6).Value = "app"
$ws.Cells.Item(57, 7).Value = "all_in_course"
$ws.Cells.Item(57, 8).Value = "THREAD_REPLY"
$ws.Cells.Item(57, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(57, 10).Value = "CourseMessageView"
$ws.Cells.Item(57, 11).NumberFormat = "@"
$ws.Cells.Item(57, 11).Value = "201"
$ws.Cells.Item(57, 12).Value = "KEEP"
$ws.Cells.Item(57, 13).Value = "返信投稿(image_32d9b1)。announcementタイプは講師以外返信不可。"
$ws.Cells.Item(57, 2).Value = "API-056"
$ws.Cells.Item(58, 3).Value = "/api/v1/messages/{messageId}"
$ws.Cells.Item(58, 4).Value = "PATCH"
$ws.Cells.Item(58, 5).Value = "メッセージ編集"
$ws.Cells.Item(58, 6).Value = "app"
$ws.Cells.Item(58, 7).Value = "owner_only"
$ws.Cells.Item(58, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(58, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(58, 10).Value = "CourseMessageView"
$ws.Cells.Item(58, 11).Value = "200, 423"
$ws.Cells.Item(58, 12).Value = "KEEP"
$ws.Cells.Item(58, 13).Value = "自分の投稿内容を修正。編集履歴を保持。423チェック対象。"
$ws.Cells.Item(58, 2).Value = "API-057"
$ws.Cells.Item(59, 3).Value = "/api/v1/messages/{messageId}"
$ws.Cells.Item(59, 4).Value = "DELETE"
$ws.Cells.Item(59, 5).Value = "メッセージ削除"
$ws.Cells.Item(59, 6).Value = "app"
$ws.Cells.Item(59, 7).Value = "owner_only, instructor"
$ws.Cells.Item(59, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(59, 9).Value = "-"
$ws.Cells.Item(59, 10).Value = "SuccessResponse"
$ws.Cells.Item(59, 11).Value = "200, 423"
$ws.Cells.Item(59, 12).Value = "KEEP"
$ws.Cells.Item(59, 13).Value = "投稿の論理削除。返信がある場合は「削除されました」と表示。"
$ws.Cells.Item(59, 2).Value = "API-058"
$ws.Cells.Item(60, 3).Value = "/api/v1/messages/{messageId}/reactions"
$ws.Cells.Item(60, 4).Value = "POST"
$ws.Cells.Item(60, 5).Value = "リアクション追加"
$ws.Cells.Item(60, 6).Value = "app"
$ws.Cells.Item(60, 7).Value = "all_in_course"
$ws.Cells.Item(60, 8).Value = "423_ON_FROZEN"
$ws.Cells.Item(60, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(60, 10).Value = "CourseMessageView"
$ws.Cells.Item(60, 11).Value = "201, 423"
$ws.Cells.Item(60, 12).Value = "KEEP"
$ws.Cells.Item(60, 13).Value = "メッセージに対する絵文字リアクションの付与。423チェック対象。"
$ws.Cells.Item(60, 2).Value = "API-059"
$ws.Cells.Item(61, 3).Value = "/api/v1/courses/{courseId}/channels/{channelId}/threads"
$ws.Cells.Item(61, 4).Value = "GET"
$ws.Cells.Item(61, 5).Value = "スレッド一覧取得（ルートメッセージのみ）"
$ws.Cells.Item(61, 6).Value = "app"
$ws.Cells.Item(61, 7).Value = "all_in_course"
$ws.Cells.Item(61, 8).Value = "-"
$ws.Cells.Item(61, 9).Value = "-"
$ws.Cells.Item(61, 10).Value = "ThreadListResponse"
$ws.Cells.Item(61, 11).NumberFormat = "@"
$ws.Cells.Item(61, 11).Value = "200"
$ws.Cells.Item(61, 12).Value = "KEEP"
$ws.Cells.Item(61, 13).Value = "threads_only対応"
$ws.Cells.Item(61, 2).Value = "API-060"
$ws.Cells.Item(62, 3).Value = "/api/v1/courses/{courseId}/channels/{channelId}/threads"
$ws.Cells.Item(62, 4).Value = "POST"
$ws.Cells.Item(62, 5).Value = "スレッド作成（ルート投稿）"
$ws.Cells.Item(62, 6).Value = "app"
$ws.Cells.Item(62, 7).Value = "all_in_course"
$ws.Cells.Item(62, 8).Value = "-"
$ws.Cells.Item(62, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(62, 10).Value = "CourseDetailView"
$ws.Cells.Item(62, 11).NumberFormat = "@"
$ws.Cells.Item(62, 11).Value = "201"
$ws.Cells.Item(62, 12).Value = "KEEP"
$ws.Cells.Item(62, 13).Value = "threadId=NULLで作成"
$ws.Cells.Item(62, 2).Value = "API-061"
$ws.Cells.Item(63, 3).Value = "/api/v1/courses/{courseId}/channels/{channelId}/threads/{threadId}/messages"
$ws.Cells.Item(63, 4).Value = "GET"
$ws.Cells.Item(63, 5).Value = "スレッド内メッセージ一覧（ルート+返信）"
$ws.Cells.Item(63, 6).Value = "app"
$ws.Cells.Item(63, 7).Value = "all_in_course"
$ws.Cells.Item(63, 8).Value = "-"
$ws.Cells.Item(63, 9).Value = "-"
$ws.Cells.Item(63, 10).Value = "MessageListResponse"
$ws.Cells.Item(63, 11).NumberFormat = "@"
$ws.Cells.Item(63, 11).Value = "200"
$ws.Cells.Item(63, 12).Value = "KEEP"
$ws.Cells.Item(63, 13).Value = "threadId指定で返信も取得"
$ws.Cells.Item(63, 2).Value = "API-062"
$ws.Cells.Item(64, 3).Value = "/api/v1/courses/{courseId}/channels/{channelId}/threads/{threadId}/messages"
$ws.Cells.Item(64, 4).Value = "POST"
$ws.Cells.Item(64, 5).Value = "スレッド返信"
$ws.Cells.Item(64, 6).Value = "app"
$ws.Cells.Item(64, 7).Value = "all_in_course"
$ws.Cells.Item(64, 8).Value = "-"
$ws.Cells.Item(64, 9).Value = "CourseMessageCreateRequest"
$ws.Cells.Item(64, 10).Value = "CourseDetailView"
$ws.Cells.Item(64, 11).NumberFormat = "@"
$ws.Cells.Item(64, 11).Value = "201"
$ws.Cells.Item(64, 12).Value = "KEEP"
$ws.Cells.Item(64, 13).Value = "threadId必須"
$ws.Cells.Item(64, 2).Value = "API-063"
$ws.Cells.Item(65, 3).Value = "/api/v1/payments/webhook"
$ws.Cells.Item(65, 4).Value = "POST"
$ws.Cells.Item(65, 5).Value = "Stripe Webhook"
$ws.Cells.Item(65, 6).Value = "app"
$ws.Cells.Item(65, 7).Value = "public(Stripe)"
$ws.Cells.Item(65, 8).Value = "Webhook専用"
$ws.Cells.Item(65, 9).Value = "GenericWriteRequest"
$ws.Cells.Item(65, 10).Value = "SuccessResponse"
$ws.Cells.Item(65, 11).NumberFormat = "@"
$ws.Cells.Item(65, 11).Value = "200"
$ws.Cells.Item(65, 12).Value = "KEEP"
$ws.Cells.Item(65, 13).Value = "外部決済完了通知を受け取り、Enrollmentを自動的にactive化する。"
$ws.Cells.Item(65, 2).Value = "API-064"
$ws.Cells.Item(66, 3).Value = "/api/v1/health"
$ws.Cells.Item(66, 4).Value = "GET"
$ws.Cells.Item(66, 5).Value = "ヘルスチェック"
$ws.Cells.Item(66, 6).Value = "app"
$ws.Cells.Item(66, 7).Value = "public"
$ws.Cells.Item(66, 8).Value = "-"
$ws.Cells.Item(66, 9).Value = "-"
$ws.Cells.Item(66, 10).Value = "GenericListResponse"
$ws.Cells.Item(66, 11).NumberFormat = "@"
$ws.Cells.Item(66, 11).Value = "200"
$ws.Cells.Item(66, 12).Value = "KEEP"
$ws.Cells.Item(66, 13).Value = "サーバーの死活監視用。DB、Redis、外部サービスの接続確認。"
$ws.Cells.Item(66, 2).Value = "API-065"
